# Auto-generated edit script: updates market-price columns (H-N) on several
# rows across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW to match refreshed
# market data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(94, 8).Value = 5733.3335
$ws.Cells.Item(94, 10).Value = 6000
$ws.Cells.Item(94, 12).Value = 6000
$ws.Cells.Item(94, 14).Value = -6902

$ws.Cells.Item(97, 8).Value = 941.6667
$ws.Cells.Item(97, 10).Value = 941.6667
$ws.Cells.Item(97, 12).Value = 2825.0001
$ws.Cells.Item(97, 14).Value = -3817.0001

$ws.Cells.Item(101, 8).Value = 564.5
$ws.Cells.Item(101, 9).Value = 516.5714
$ws.Cells.Item(101, 10).Value = 900
$ws.Cells.Item(101, 11).Value = 1549.7142
$ws.Cells.Item(101, 12).Value = 2700
$ws.Cells.Item(101, 13).Value = 72.28579999999988
$ws.Cells.Item(101, 14).Value = -5944

$ws.Cells.Item(112, 8).Value = 1441.7872
$ws.Cells.Item(112, 9).Value = 921.17645
$ws.Cells.Item(112, 10).Value = 1736.8
$ws.Cells.Item(112, 11).Value = 2763.52935
$ws.Cells.Item(112, 12).Value = 5210.4
$ws.Cells.Item(112, 13).Value = -1655.52935
$ws.Cells.Item(112, 14).Value = -7426.4

$ws.Cells.Item(129, 8).Value = 799.26
$ws.Cells.Item(129, 9).Value = 598.6316
$ws.Cells.Item(129, 10).Value = 922.2258
$ws.Cells.Item(129, 11).Value = 1795.8948
$ws.Cells.Item(129, 12).Value = 2766.6774
$ws.Cells.Item(129, 13).Value = 3204.1052
$ws.Cells.Item(129, 14).Value = -12766.6774

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 1851.25
$ws.Cells.Item(63, 9).Value = 1635
$ws.Cells.Item(63, 10).Value = 2500
$ws.Cells.Item(63, 11).Value = 1635
$ws.Cells.Item(63, 12).Value = 2500
$ws.Cells.Item(63, 13).Value = -949
$ws.Cells.Item(63, 14).Value = -3872

$ws.Cells.Item(66, 8).Value = 1851.25
$ws.Cells.Item(66, 9).Value = 1635
$ws.Cells.Item(66, 10).Value = 2500
$ws.Cells.Item(66, 11).Value = 8175
$ws.Cells.Item(66, 12).Value = 12500
$ws.Cells.Item(66, 13).Value = -4743
$ws.Cells.Item(66, 14).Value = -19364

$ws.Cells.Item(88, 8).Value = 2690
$ws.Cells.Item(88, 9).Value = 3434.3333
$ws.Cells.Item(88, 10).Value = 2487
$ws.Cells.Item(88, 11).Value = 3434.3333
$ws.Cells.Item(88, 12).Value = 2487
$ws.Cells.Item(88, 13).Value = -3028.3333
$ws.Cells.Item(88, 14).Value = -3299

$ws.Cells.Item(91, 8).Value = 2690
$ws.Cells.Item(91, 9).Value = 3434.3333
$ws.Cells.Item(91, 10).Value = 2487
$ws.Cells.Item(91, 11).Value = 3434.3333
$ws.Cells.Item(91, 12).Value = 2487
$ws.Cells.Item(91, 13).Value = -2030.3333
$ws.Cells.Item(91, 14).Value = -5295

$ws.Cells.Item(97, 8).Value = 633.24
$ws.Cells.Item(97, 9).Value = 465.2353
$ws.Cells.Item(97, 10).Value = 990.25
$ws.Cells.Item(97, 11).Value = 465.2353
$ws.Cells.Item(97, 12).Value = 990.25
$ws.Cells.Item(97, 13).Value = 30.7647
$ws.Cells.Item(97, 14).Value = -1982.25

$ws.Cells.Item(132, 8).Value = 9222.786
$ws.Cells.Item(132, 9).Value = 8241.6
$ws.Cells.Item(132, 10).Value = 10354.923
$ws.Cells.Item(132, 11).Value = 24724.8
$ws.Cells.Item(132, 12).Value = 31064.769
$ws.Cells.Item(132, 13).Value = -22194.8
$ws.Cells.Item(132, 14).Value = -36124.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 401847.2
$ws.Cells.Item(86, 9).Value = 1565.5555
$ws.Cells.Item(86, 10).Value = 1431142.9
$ws.Cells.Item(86, 11).Value = 1565.5555
$ws.Cells.Item(86, 12).Value = 1431142.9
$ws.Cells.Item(86, 13).Value = -442.5554999999999
$ws.Cells.Item(86, 14).Value = -1433388.9

$ws.Cells.Item(89, 8).Value = 401847.2
$ws.Cells.Item(89, 9).Value = 1565.5555
$ws.Cells.Item(89, 10).Value = 1431142.9
$ws.Cells.Item(89, 11).Value = 7827.7775
$ws.Cells.Item(89, 12).Value = 7155714.5
$ws.Cells.Item(89, 13).Value = -2211.7775
$ws.Cells.Item(89, 14).Value = -7166946.5

$ws.Cells.Item(94, 8).Value = 1058.76
$ws.Cells.Item(94, 9).Value = 829.64703
$ws.Cells.Item(94, 10).Value = 1545.625
$ws.Cells.Item(94, 11).Value = 829.64703
$ws.Cells.Item(94, 12).Value = 1545.625
$ws.Cells.Item(94, 13).Value = -378.64703
$ws.Cells.Item(94, 14).Value = -2447.625

$ws.Cells.Item(99, 8).Value = 1495.5385
$ws.Cells.Item(99, 9).Value = 835
$ws.Cells.Item(99, 10).Value = 2061.7144
$ws.Cells.Item(99, 11).Value = 835
$ws.Cells.Item(99, 12).Value = 2061.7144
$ws.Cells.Item(99, 13).Value = 663
$ws.Cells.Item(99, 14).Value = -5057.7144

$ws.Cells.Item(105, 8).Value = 1678.4375
$ws.Cells.Item(105, 9).Value = 1358.5454
$ws.Cells.Item(105, 11).Value = 1358.5454
$ws.Cells.Item(105, 13).Value = 388.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()

$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 893.67444
$ws.Cells.Item(131, 10).Value = 958.25
$ws.Cells.Item(131, 12).Value = 2874.75
$ws.Cells.Item(131, 14).Value = -12954.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1033.3334
$ws.Cells.Item(122, 9).Value = 966.6667
$ws.Cells.Item(122, 10).Value = 1166.6666
$ws.Cells.Item(122, 11).Value = 2900.0001
$ws.Cells.Item(122, 12).Value = 3499.9998
$ws.Cells.Item(122, 13).Value = -450.0001000000002
$ws.Cells.Item(122, 14).Value = -8399.9998

$ws.Cells.Item(132, 8).Value = 15607.8
$ws.Cells.Item(132, 9).Value = 41508.668
$ws.Cells.Item(132, 10).Value = 4507.4287
$ws.Cells.Item(132, 11).Value = 124526.004
$ws.Cells.Item(132, 12).Value = 13522.2861
$ws.Cells.Item(132, 13).Value = -121996.004
$ws.Cells.Item(132, 14).Value = -18582.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1678.0714
$ws.Cells.Item(7, 9).Value = 1451.875
$ws.Cells.Item(7, 10).Value = 1979.6666
$ws.Cells.Item(7, 11).Value = 1451.875
$ws.Cells.Item(7, 12).Value = 1979.6666
$ws.Cells.Item(7, 13).Value = -1339.875
$ws.Cells.Item(7, 14).Value = -2203.6666

$ws.Cells.Item(40, 8).Value = 2469.8333
$ws.Cells.Item(40, 9).Value = 1847
$ws.Cells.Item(40, 11).Value = 1847
$ws.Cells.Item(40, 13).Value = -1711

$ws.Cells.Item(55, 8).Value = 269.85715
$ws.Cells.Item(55, 9).Value = 147
$ws.Cells.Item(55, 10).Value = 433.66666
$ws.Cells.Item(55, 11).Value = 147
$ws.Cells.Item(55, 12).Value = 433.66666
$ws.Cells.Item(55, 13).Value = 26
$ws.Cells.Item(55, 14).Value = -779.66666

$ws.Cells.Item(122, 8).Value = 22003.373
$ws.Cells.Item(122, 9).Value = 32546.697
$ws.Cells.Item(122, 10).Value = 2673.9443
$ws.Cells.Item(122, 11).Value = 97640.091
$ws.Cells.Item(122, 12).Value = 8021.8329
$ws.Cells.Item(122, 13).Value = -95190.091
$ws.Cells.Item(122, 14).Value = -12921.8329

$ws.Cells.Item(126, 8).Value = 1678.0714
$ws.Cells.Item(126, 9).Value = 1451.875
$ws.Cells.Item(126, 10).Value = 1979.6666
$ws.Cells.Item(126, 11).Value = 4355.625
$ws.Cells.Item(126, 12).Value = 5938.9998
$ws.Cells.Item(126, 13).Value = -1885.625
$ws.Cells.Item(126, 14).Value = -10878.9998
